# Updating the Staging testdata
$wb = $excel.ActiveWorkbook

$wsNew = $wb.Worksheets.Item("NewImportLogic")
$wsOld = $wb.Worksheets.Item("OldImportLogic")

# Update the cell text on the NewImportLogic sheet (remove spaces around dash)
$wsNew.Range("H2").Value = "ExcelReport-QOL_and_ECON-UtilityOutcome-Economic-"

# Set selection on OldImportLogic sheet (stays at H3, unchanged)
$wsOld.Activate()
$wsOld.Range("H3").Select()

# Set selection on NewImportLogic and make it the active tab
$wsNew.Activate()
$wsNew.Range("H3").Select()

$wb.Save()
